$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Instal days by Model")

# Insert two new rows above the former row 29 ("SS25") for the new
# RPC-PH / RPC-OU items, shifting the rest of the table down.
$ws1.Rows("29:30").Insert()

$ws1.Range("A29").Value = "RPC-PH"
$ws1.Range("B29").Value = 4
$ws1.Range("C29").Value = 2
$ws1.Range("D29").Value = $true
$ws1.Range("E29").Value = $true

$ws1.Range("A30").Value = "RPC-OU"
$ws1.Range("B30").Value = 4
$ws1.Range("C30").Value = 2
$ws1.Range("D30").Value = $true
$ws1.Range("E30").Value = $true

# Grow Table1 (and its AutoFilter) to cover the two new rows.
$lo = $ws1.ListObjects.Item("Table1")
$lo.Resize($ws1.Range("A1:E44"))

# Keep the hidden _FilterDatabase name in sync with the new table size.
$dn = $wb.Names.Item("_xlnm._FilterDatabase")
$dn.RefersTo = "='Instal days by Model'!`$A`$1:`$C`$44"

# Make "Instal days by Model" the active sheet/selection (was "Service
# Rates" before), with A31 selected.
$ws1.Activate()
$ws1.Range("A31").Select()
